# Export with no is_pref and no lev distance
# Re-map each row's id/speaker_variant to its canonical (self-matching) tag
# and clear all "is_prefered" (column D) marks.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = "#tiel"
$ws.Range("C2").Value = "Tiel"
$ws.Range("D2").Value = ""

$ws.Range("B3").Value = "#alet"
$ws.Range("C3").Value = "Alet"
$ws.Range("D3").Value = ""

$ws.Range("B4").Value = "#brect"
$ws.Range("C4").Value = "Brect"
$ws.Range("D4").Value = ""

$ws.Range("B5").Value = "#sol"
$ws.Range("C5").Value = "Sol"
$ws.Range("D5").Value = ""

$ws.Range("B6").Value = "#haripon"
$ws.Range("C6").Value = "Haripon"
$ws.Range("D6").Value = ""

$ws.Range("B7").Value = "#kaat"
$ws.Range("C7").Value = "Kaat"
$ws.Range("D7").Value = ""

$ws.Range("D8").Value = ""

$ws.Range("B9").Value = "#jodelet"
$ws.Range("C9").Value = "Jodelet"
$ws.Range("D9").Value = ""

$ws.Range("B10").Value = "#marotte"
$ws.Range("C10").Value = "Marotte"
$ws.Range("D10").Value = ""

$ws.Range("B11").Value = "#alardus"
$ws.Range("C11").Value = "Alardus"
$ws.Range("D11").Value = ""

$ws.Range("B12").Value = "#ferdinand"
$ws.Range("C12").Value = "Ferdinand"
$ws.Range("D12").Value = ""

$ws.Range("B13").Value = "#frederik"
$ws.Range("C13").Value = "Frederik"
$ws.Range("D13").Value = ""

$ws.Range("D14").Value = ""

$ws.Range("B15").Value = "#boon"
$ws.Range("C15").Value = "Boon"

$ws.Range("B16").Value = "#roel"
$ws.Range("C16").Value = "Roel"

$ws.Range("C17").Value = "Piet"

$ws.Range("B18").Value = "#dina"
$ws.Range("C18").Value = "Dina"
